$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 93, shifting existing rows 93:164 down to 94:165
$ws.Rows.Item(93).Insert()

# Fill in the newly inserted row 93 with the new data record
$ws.Cells.Item(93, 1).Value  = 4
$ws.Cells.Item(93, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(93, 3).Value  = "Los Lagos"
$ws.Cells.Item(93, 4).Value  = 44957
$ws.Cells.Item(93, 5).Value  = 10
$ws.Cells.Item(93, 6).Value  = 100112052
$ws.Cells.Item(93, 7).Value  = "Albahaca"
$ws.Cells.Item(93, 8).Value  = "Sin especificar"
$ws.Cells.Item(93, 9).Value  = "Primera"
$ws.Cells.Item(93, 10).Value = 120
$ws.Cells.Item(93, 11).Value = 6000
$ws.Cells.Item(93, 12).Value = 6000
$ws.Cells.Item(93, 13).Value = 6000
$ws.Cells.Item(93, 14).Value = "$/docena de matas"
$ws.Cells.Item(93, 15).Value = "Región Metropolitana"
$ws.Cells.Item(93, 16).Value = 1000
$ws.Cells.Item(93, 17).Value = 6
$ws.Cells.Item(93, 18).Value = "Hortaliza"

# Match the date-number-format style used by the other rows in column D
$ws.Cells.Item(93, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
